$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-08-25"
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = 57.43000030517578
$ws.Range("C10").Value = 686.7999877929688
$ws.Range("D10").Value = 319.25
